# Meeting20150425/Lack of manpower.pptx - "Added availability of JDW."
#
# 1) Re-style every table (slides 2-8, one table each) with a new
#    table-style GUID (as if a different style was picked from the
#    Table Styles gallery for each table).
# 2) Append ", not before July" (24pt) after "JDW" in the title of
#    slide 7.

$p = $ppt.ActivePresentation

# --- 1. Table style swaps -------------------------------------------------
$tableStyleUpdates = @(
    @{ Slide = 2; NewStyle = "{0BC1FBB6-F04C-420D-B3E1-8E79F56C9579}" },
    @{ Slide = 3; NewStyle = "{90AD7B7B-C225-4A48-ACD4-69576CDBA7E0}" },
    @{ Slide = 4; NewStyle = "{70D44EBE-3112-42A9-B50B-D39B157FB8B3}" },
    @{ Slide = 5; NewStyle = "{1A805F8A-0E14-44CD-86D9-8F55B72DDFF6}" },
    @{ Slide = 6; NewStyle = "{AAA811F6-21D2-44E7-B18C-B693A088D4EE}" },
    @{ Slide = 7; NewStyle = "{0B4E03C1-9374-4960-AE57-F2553D14477F}" },
    @{ Slide = 8; NewStyle = "{2000785D-C8BD-4395-81AC-27520309F939}" }
)

foreach ($update in $tableStyleUpdates) {
    $slide = $p.Slides.Item($update.Slide)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($update.NewStyle)
        }
    }
}

# --- 2. Title text edit on slide 7 ----------------------------------------
$slide7 = $p.Slides.Item(7)
$titleShape = $slide7.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$originalLength = $titleRange.Length
$titleRange.InsertAfter(", not before July") | Out-Null

$newRunRange = $titleRange.Characters($originalLength + 1, 18)
$newRunRange.Font.Size = 24
